$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply 2020-07-28 data updates to nombre_aides (col C) and montant_total (col D).
# NumberFormat is forced to text ("@") before the write so the numeric-looking
# strings are kept as text (matching the source inlineStr cells) instead of being
# auto-converted to numbers; Style is reset back to Normal afterwards so no stray
# cell formatting is introduced.
$c = $ws.Cells.Item(8, 3)
$c.NumberFormat = "@"
$c.Value = "615"
$c.Style = "Normal"
$d = $ws.Cells.Item(8, 4)
$d.NumberFormat = "@"
$d.Value = "1893584.36"
$d.Style = "Normal"
$c = $ws.Cells.Item(48, 3)
$c.NumberFormat = "@"
$c.Value = "22"
$c.Style = "Normal"
$d = $ws.Cells.Item(48, 4)
$d.NumberFormat = "@"
$d.Value = "63500.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(50, 3)
$c.NumberFormat = "@"
$c.Value = "105"
$c.Style = "Normal"
$d = $ws.Cells.Item(50, 4)
$d.NumberFormat = "@"
$d.Value = "271468.33"
$d.Style = "Normal"
$c = $ws.Cells.Item(52, 3)
$c.NumberFormat = "@"
$c.Value = "233"
$c.Style = "Normal"
$d = $ws.Cells.Item(52, 4)
$d.NumberFormat = "@"
$d.Value = "752162.59"
$d.Style = "Normal"
$c = $ws.Cells.Item(56, 3)
$c.NumberFormat = "@"
$c.Value = "33"
$c.Style = "Normal"
$d = $ws.Cells.Item(56, 4)
$d.NumberFormat = "@"
$d.Value = "70500.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(60, 3)
$c.NumberFormat = "@"
$c.Value = "34"
$c.Style = "Normal"
$d = $ws.Cells.Item(60, 4)
$d.NumberFormat = "@"
$d.Value = "124456.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(76, 3)
$c.NumberFormat = "@"
$c.Value = "75"
$c.Style = "Normal"
$d = $ws.Cells.Item(76, 4)
$d.NumberFormat = "@"
$d.Value = "199487.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(77, 3)
$c.NumberFormat = "@"
$c.Value = "171"
$c.Style = "Normal"
$d = $ws.Cells.Item(77, 4)
$d.NumberFormat = "@"
$d.Value = "446402.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(78, 3)
$c.NumberFormat = "@"
$c.Value = "21"
$c.Style = "Normal"
$d = $ws.Cells.Item(78, 4)
$d.NumberFormat = "@"
$d.Value = "63991.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(79, 3)
$c.NumberFormat = "@"
$c.Value = "368"
$c.Style = "Normal"
$d = $ws.Cells.Item(79, 4)
$d.NumberFormat = "@"
$d.Value = "1256755.04"
$d.Style = "Normal"
$c = $ws.Cells.Item(82, 3)
$c.NumberFormat = "@"
$c.Value = "12"
$c.Style = "Normal"
$d = $ws.Cells.Item(82, 4)
$d.NumberFormat = "@"
$d.Value = "36500.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(83, 3)
$c.NumberFormat = "@"
$c.Value = "56"
$c.Style = "Normal"
$d = $ws.Cells.Item(83, 4)
$d.NumberFormat = "@"
$d.Value = "171200.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(85, 3)
$c.NumberFormat = "@"
$c.Value = "37"
$c.Style = "Normal"
$d = $ws.Cells.Item(85, 4)
$d.NumberFormat = "@"
$d.Value = "87500.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(87, 3)
$c.NumberFormat = "@"
$c.Value = "48"
$c.Style = "Normal"
$d = $ws.Cells.Item(87, 4)
$d.NumberFormat = "@"
$d.Value = "182846.08"
$d.Style = "Normal"
$c = $ws.Cells.Item(104, 3)
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"
$d = $ws.Cells.Item(104, 4)
$d.NumberFormat = "@"
$d.Value = "15351.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(105, 3)
$c.NumberFormat = "@"
$c.Value = "10"
$c.Style = "Normal"
$d = $ws.Cells.Item(105, 4)
$d.NumberFormat = "@"
$d.Value = "30709.84"
$d.Style = "Normal"
$c = $ws.Cells.Item(106, 3)
$c.NumberFormat = "@"
$c.Value = "31"
$c.Style = "Normal"
$d = $ws.Cells.Item(106, 4)
$d.NumberFormat = "@"
$d.Value = "85860.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(107, 3)
$c.NumberFormat = "@"
$c.Value = "17"
$c.Style = "Normal"
$d = $ws.Cells.Item(107, 4)
$d.NumberFormat = "@"
$d.Value = "58850.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(109, 3)
$c.NumberFormat = "@"
$c.Value = "47"
$c.Style = "Normal"
$d = $ws.Cells.Item(109, 4)
$d.NumberFormat = "@"
$d.Value = "286558.15"
$d.Style = "Normal"
$c = $ws.Cells.Item(111, 3)
$c.NumberFormat = "@"
$c.Value = "13"
$c.Style = "Normal"
$d = $ws.Cells.Item(111, 4)
$d.NumberFormat = "@"
$d.Value = "46267.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(112, 3)
$c.NumberFormat = "@"
$c.Value = "15"
$c.Style = "Normal"
$d = $ws.Cells.Item(112, 4)
$d.NumberFormat = "@"
$d.Value = "41895.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(115, 3)
$c.NumberFormat = "@"
$c.Value = "15"
$c.Style = "Normal"
$d = $ws.Cells.Item(115, 4)
$d.NumberFormat = "@"
$d.Value = "56556.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(116, 3)
$c.NumberFormat = "@"
$c.Value = "10"
$c.Style = "Normal"
$d = $ws.Cells.Item(116, 4)
$d.NumberFormat = "@"
$d.Value = "20000.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(117, 3)
$c.NumberFormat = "@"
$c.Value = "27"
$c.Style = "Normal"
$d = $ws.Cells.Item(117, 4)
$d.NumberFormat = "@"
$d.Value = "90500.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(118, 3)
$c.NumberFormat = "@"
$c.Value = "58"
$c.Style = "Normal"
$d = $ws.Cells.Item(118, 4)
$d.NumberFormat = "@"
$d.Value = "161377.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(119, 3)
$c.NumberFormat = "@"
$c.Value = "207"
$c.Style = "Normal"
$d = $ws.Cells.Item(119, 4)
$d.NumberFormat = "@"
$d.Value = "572500.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(120, 3)
$c.NumberFormat = "@"
$c.Value = "46"
$c.Style = "Normal"
$d = $ws.Cells.Item(120, 4)
$d.NumberFormat = "@"
$d.Value = "132547.58"
$d.Style = "Normal"
$c = $ws.Cells.Item(121, 3)
$c.NumberFormat = "@"
$c.Value = "381"
$c.Style = "Normal"
$d = $ws.Cells.Item(121, 4)
$d.NumberFormat = "@"
$d.Value = "1485612.45"
$d.Style = "Normal"
$c = $ws.Cells.Item(126, 3)
$c.NumberFormat = "@"
$c.Value = "33"
$c.Style = "Normal"
$d = $ws.Cells.Item(126, 4)
$d.NumberFormat = "@"
$d.Value = "138579.76"
$d.Style = "Normal"
$c = $ws.Cells.Item(130, 3)
$c.NumberFormat = "@"
$c.Value = "103"
$c.Style = "Normal"
$d = $ws.Cells.Item(130, 4)
$d.NumberFormat = "@"
$d.Value = "257519.44"
$d.Style = "Normal"
$c = $ws.Cells.Item(182, 3)
$c.NumberFormat = "@"
$c.Value = "328"
$c.Style = "Normal"
$d = $ws.Cells.Item(182, 4)
$d.NumberFormat = "@"
$d.Value = "873788.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(184, 3)
$c.NumberFormat = "@"
$c.Value = "562"
$c.Style = "Normal"
$d = $ws.Cells.Item(184, 4)
$d.NumberFormat = "@"
$d.Value = "1940016.91"
$d.Style = "Normal"
$c = $ws.Cells.Item(190, 3)
$c.NumberFormat = "@"
$c.Value = "66"
$c.Style = "Normal"
$d = $ws.Cells.Item(190, 4)
$d.NumberFormat = "@"
$d.Value = "148000.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(192, 3)
$c.NumberFormat = "@"
$c.Value = "95"
$c.Style = "Normal"
$d = $ws.Cells.Item(192, 4)
$d.NumberFormat = "@"
$d.Value = "368012.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(193, 3)
$c.NumberFormat = "@"
$c.Value = "124"
$c.Style = "Normal"
$d = $ws.Cells.Item(193, 4)
$d.NumberFormat = "@"
$d.Value = "267196.77"
$d.Style = "Normal"
$c = $ws.Cells.Item(209, 3)
$c.NumberFormat = "@"
$c.Value = "15"
$c.Style = "Normal"
$d = $ws.Cells.Item(209, 4)
$d.NumberFormat = "@"
$d.Value = "44000.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(210, 3)
$c.NumberFormat = "@"
$c.Value = "22"
$c.Style = "Normal"
$d = $ws.Cells.Item(210, 4)
$d.NumberFormat = "@"
$d.Value = "60523.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(211, 3)
$c.NumberFormat = "@"
$c.Value = "56"
$c.Style = "Normal"
$d = $ws.Cells.Item(211, 4)
$d.NumberFormat = "@"
$d.Value = "149500.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(212, 3)
$c.NumberFormat = "@"
$c.Value = "153"
$c.Style = "Normal"
$d = $ws.Cells.Item(212, 4)
$d.NumberFormat = "@"
$d.Value = "421905.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(213, 3)
$c.NumberFormat = "@"
$c.Value = "14"
$c.Style = "Normal"
$d = $ws.Cells.Item(213, 4)
$d.NumberFormat = "@"
$d.Value = "36000.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(214, 3)
$c.NumberFormat = "@"
$c.Value = "304"
$c.Style = "Normal"
$d = $ws.Cells.Item(214, 4)
$d.NumberFormat = "@"
$d.Value = "993845.50"
$d.Style = "Normal"
$c = $ws.Cells.Item(218, 3)
$c.NumberFormat = "@"
$c.Value = "72"
$c.Style = "Normal"
$d = $ws.Cells.Item(218, 4)
$d.NumberFormat = "@"
$d.Value = "203687.09"
$d.Style = "Normal"
$c = $ws.Cells.Item(219, 3)
$c.NumberFormat = "@"
$c.Value = "24"
$c.Style = "Normal"
$d = $ws.Cells.Item(219, 4)
$d.NumberFormat = "@"
$d.Value = "76587.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(222, 3)
$c.NumberFormat = "@"
$c.Value = "46"
$c.Style = "Normal"
$d = $ws.Cells.Item(222, 4)
$d.NumberFormat = "@"
$d.Value = "141270.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(223, 3)
$c.NumberFormat = "@"
$c.Value = "70"
$c.Style = "Normal"
$d = $ws.Cells.Item(223, 4)
$d.NumberFormat = "@"
$d.Value = "147500.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(225, 3)
$c.NumberFormat = "@"
$c.Value = "72"
$c.Style = "Normal"
$d = $ws.Cells.Item(225, 4)
$d.NumberFormat = "@"
$d.Value = "192538.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(226, 3)
$c.NumberFormat = "@"
$c.Value = "133"
$c.Style = "Normal"
$d = $ws.Cells.Item(226, 4)
$d.NumberFormat = "@"
$d.Value = "347200.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(227, 3)
$c.NumberFormat = "@"
$c.Value = "442"
$c.Style = "Normal"
$d = $ws.Cells.Item(227, 4)
$d.NumberFormat = "@"
$d.Value = "1142583.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(228, 3)
$c.NumberFormat = "@"
$c.Value = "69"
$c.Style = "Normal"
$d = $ws.Cells.Item(228, 4)
$d.NumberFormat = "@"
$d.Value = "190788.64"
$d.Style = "Normal"
$c = $ws.Cells.Item(229, 3)
$c.NumberFormat = "@"
$c.Value = "842"
$c.Style = "Normal"
$d = $ws.Cells.Item(229, 4)
$d.NumberFormat = "@"
$d.Value = "2642473.62"
$d.Style = "Normal"
$c = $ws.Cells.Item(234, 3)
$c.NumberFormat = "@"
$c.Value = "106"
$c.Style = "Normal"
$d = $ws.Cells.Item(234, 4)
$d.NumberFormat = "@"
$d.Value = "300476.00"
$d.Style = "Normal"
$c = $ws.Cells.Item(236, 3)
$c.NumberFormat = "@"
$c.Value = "28"
$c.Style = "Normal"
$d = $ws.Cells.Item(236, 4)
$d.NumberFormat = "@"
$d.Value = "66500.00"
$d.Style = "Normal"
